# Added Password Encrypt and Decrypt logic
# Replace plaintext credentials stored in the TGL test-data sheet with their
# encrypted equivalents.
#
#   "Rockstar1" (password column E on every row) -> "SyiqZDyGjHvML1WKNrXshA=="
#   "password"  (password column G on rows 2 & 5) -> "WZrRgv7ejKOtV0KoRQsURQ=="
#
# The order in which the *first* occurrence of each new value is written
# controls where it lands in the shared-strings table, so write the column E
# values (encrypted "Rockstar1") before the column G values (encrypted
# "password") to match the original authoring order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$encryptedRockstar = "SyiqZDyGjHvML1WKNrXshA=="
$encryptedPassword = "WZrRgv7ejKOtV0KoRQsURQ=="

$ws.Range("E2").Value = $encryptedRockstar
$ws.Range("E3").Value = $encryptedRockstar
$ws.Range("E4").Value = $encryptedRockstar
$ws.Range("E5").Value = $encryptedRockstar

$ws.Range("G3").Value = $encryptedRockstar
$ws.Range("G4").Value = $encryptedRockstar

$ws.Range("G2").Value = $encryptedPassword
$ws.Range("G5").Value = $encryptedPassword

# Move the active selection to G2 (matches the saved cursor position after
# the edit) and drop the stale F1 "frozen" top-left scroll position.
$ws.Range("G2").Select()
